$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.554.84"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.362.74"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'256.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'665.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.41%  "
$ws.Range("D7").Value = "'1.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.61%  "
$ws.Range("D8").Value = "'0.471"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.30%  "
$ws.Range("D9").Value = "'1.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +22.15%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "3.359.87"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  +8.63%  "
$ws.Range("D13").Value = "'42.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.44%  "
$ws.Range("D14").Value = "'0.0000273"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.60%  "
$ws.Range("D15").Value = "98.892.01"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "3.985.76"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "'5.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +27.81%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.352.54"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'16.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.61%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'532.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.45%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "'3.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "'10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.65%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "'0.438"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +50.46%  "
$ws.Range("D26").Value = "'102.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.14%  "
$ws.Range("D27").Value = "'6.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.14%  "
$ws.Range("E28").Value = "  +5.92%  "
$ws.Range("D29").Value = "3.541.07"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "'0.151"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.28%  "
$ws.Range("D32").Value = "'11.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.67%  "
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("D36").Value = "'0.542"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.74%  "
$ws.Range("D37").Value = "'7.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("E38").Value = "  +9.04%  "
$ws.Range("D39").Value = "'0.158"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").Value = "'524.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("B41").Value = "MantraDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D41").Value = "'3.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.67%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'24.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "'0.0436"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +34.26%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "'0.828"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +6.53%  "
$ws.Range("D49").Value = "'1.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.33%  "
$ws.Range("D50").Value = "'5.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.12%  "
$ws.Range("D51").Value = "'51.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.49%  "
